$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 17 (shifts existing rows down)
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row with the Glauciane record (correct balance 590.69).
# Force text format first so the leading zeros in the account number survive
# (Excel would otherwise parse "005981575" as a number), then restore the
# cell's formatting back to the plain default style used by every other data
# row (copy the General-format from a neighboring cell via Paste Special).
$ws.Cells.Item(17, 1).NumberFormat = "@"
$ws.Cells.Item(17, 1).Value = "005981575"
$ws.Cells.Item(17, 2).Value = "Glauciane"
$ws.Cells.Item(17, 3).Value = 590.69

$ws.Cells.Item(16, 1).Copy()
$ws.Cells.Item(17, 1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# The old Glauciane row (previously row 69, now shifted to row 70) is now a duplicate
# with the stale balance of 90.69 - remove it entirely.
$ws.Rows.Item(70).Delete()
